$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B" = 0.99999896985964509
    "C" = 0.99902446656693789
    "D" = 0.99999444230408396
    "E" = 0.99999420209957157
    "F" = 0.9999953044873553
    "G" = 0.0000009615908033740045
    "H" = 0.0009106176387947784
    "I" = 0.0000002763649468499002
    "J" = 0.000004230617425971275
    "K" = 0.000002253491186410587
    "L" = 0.00005130350677945762
    "M" = 0.000980607364531801
    "N" = 1.000024723368518
    "O" = 0.001022353850177004
    "P" = 77.70935367418030637
    "Q" = 108.18124929588529426
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
